$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from column N (rows 3, 5-25) into column O so the new
# column inherits the same cell styles as its neighbour, then fill in
# the 2021 data values.
$ws.Range("N3").Copy()
$ws.Range("O3").PasteSpecial(-4122)

$ws.Range("N5:N25").Copy()
$ws.Range("O5").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

$ws.Range("O3").Value = 2021

$ws.Range("O5").Value = 2148.2
$ws.Range("O6").Value = 109.5
$ws.Range("O7").Value = 210.1
$ws.Range("O8").Value = 196
$ws.Range("O9").Value = 209
$ws.Range("O10").Value = 300.2
$ws.Range("O11").Value = 302.9
$ws.Range("O12").Value = 786
$ws.Range("O13").Value = 27.7
$ws.Range("O14").Value = 6.8

$ws.Range("O16").Value = 26.9
$ws.Range("O17").Value = 15.9
$ws.Range("O18").Value = 21.7
$ws.Range("O19").Value = 29.9
$ws.Range("O20").Value = 30.2
$ws.Range("O21").Value = 24
$ws.Range("O22").Value = 31.6
$ws.Range("O23").Value = 30.3
$ws.Range("O24").Value = 20.7
$ws.Range("O25").Value = 12

# Move the selection as recorded in the saved workbook.
$ws.Range("Q20").Select()
